# (hotfix) record messages into Excel file
#
# The messenger logger had written the two buffered messages to the sheet
# with their timestamps/text swapped. Re-record row 2 and row 3 with the
# correct pairing of "when" (column A) and "what" (column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the text currently sitting in D2/D3 before either is overwritten.
$origD2 = $ws.Range("D2").Text
$origD3 = $ws.Range("D3").Text

# Row 2 should hold the later message...
$ws.Range("A2").Value = 45051.61481586432
$ws.Range("D2").Value = $origD3

# ...and row 3 the earlier one, recorded with its precise timestamp.
$ws.Range("A3").Value = 45051.61453230146
$ws.Range("D3").Value = $origD2

# Both timestamp cells share the same date/time number format.
$ws.Range("A3").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
